# Daily "advance the clock by one day" update for the 剩余(remaining)/开始时间(start date)
# tracking columns (E / F) on the active worksheet.
#
# For every data row (row 2 .. last used row):
#   - D = 总天 (total days for this cycle)
#   - E = 剩余 (days remaining in the current cycle)
#   - F = 开始时间 (cycle start date, stored as an 8-digit YYYYMMDD number)
#
# Each day that passes:
#   - if the remaining-days counter has already hit 1 (i.e. the cycle is
#     finishing today), the cycle rolls over: a new cycle starts on the
#     day the old cycle was due (F + D days), with the counter reset to
#     the full D days remaining.
#   - otherwise the counter simply ticks down by one, the start date is
#     unchanged.
#
# Rows whose start date cannot be parsed as a real calendar date (data
# entry typos, e.g. a 9-digit value) are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    # NOTE: `.Value` misbehaves for reads in this host (returns the
    # property descriptor instead of the cell's content) - use `.Value2`
    # for every read. Writes work fine either way; `.Value2` is used
    # throughout for consistency.
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $totalDays = [int]$dVal
    $remaining = [int]$eVal
    $startStr = [string][int64]$fVal

    if ($startStr.Length -ne 8) {
        continue
    }

    $year = [int]$startStr.Substring(0, 4)
    $month = [int]$startStr.Substring(4, 2)
    $day = [int]$startStr.Substring(6, 2)

    $startDate = $null
    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        $startDate = $null
    }

    if ($startDate -eq $null) {
        # Unparseable / malformed start date (e.g. typo) - leave row alone.
        continue
    }

    $dueDate = $startDate.AddDays($totalDays)

    if ($remaining -le 1) {
        # Cycle finished - roll over to a fresh cycle starting on the due date.
        $newStart = $dueDate
        $eCell.Value2 = $totalDays
        $fCell.Value2 = [int]($newStart.ToString("yyyyMMdd"))
    } else {
        # One more day has elapsed in the current cycle.
        $eCell.Value2 = $remaining - 1
    }
}
